# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.136.21"
$ws.Range("E2").Value = "  +0.11%  "

$ws.Range("D3").Value = "2.050.14"
$ws.Range("E3").Value = "  -1.26%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.10"
$ws.Range("E5").Value = "  -1.97%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.664"
$ws.Range("E6").Value = "  -1.29%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.26"
$ws.Range("E7").Value = "  -1.55%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.382"
$ws.Range("E9").Value = "  -1.94%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0776"
$ws.Range("E10").Value = "  -2.36%  "

$ws.Range("E11").Value = "  +0.13%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.93"
$ws.Range("E12").Value = "  +0.00%  "

$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.849"
$ws.Range("E13").Value = "  +4.06%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.348.58"
$ws.Range("E14").Value = "  -1.34%  "

$ws.Range("E15").Value = "  +3.20%  "

$ws.Range("D16").Value = "2.049.13"
$ws.Range("E16").Value = "  -1.31%  "

$ws.Range("E17").Value = "  +16.05%  "

$ws.Range("D18").Value = "37.130.28"
$ws.Range("E18").Value = "  +0.18%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.95"
$ws.Range("E19").Value = "  +0.79%  "

$ws.Range("E20").Value = "  -3.34%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.34"
$ws.Range("E21").Value = "  -1.91%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "237.17"
$ws.Range("E22").Value = "  -0.88%  "

$ws.Range("E24").Value = "  +1.97%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.49"
$ws.Range("E25").Value = "  +2.13%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.53"
$ws.Range("E26").Value = "  -0.10%  "

$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.17"
$ws.Range("E27").Value = "  -5.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.05"
$ws.Range("E28").Value = "  -1.27%  "

$ws.Range("E29").Value = "  -1.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.80"
$ws.Range("E30").Value = "  -0.58%  "

$ws.Range("E31").Value = "  -0.69%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0617"
$ws.Range("E32").Value = "  -2.91%  "

$ws.Range("E33").Value = "  +1.07%  "

$ws.Range("E34").Value = "  -2.13%  "

$ws.Range("E35").Value = "  -0.12%  "

$ws.Range("E36").Value = "  -2.27%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.76"
$ws.Range("E37").Value = "  -0.45%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.27"
$ws.Range("E38").Value = "  +15.92%  "

$ws.Range("E39").Value = "  -2.13%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.18"
$ws.Range("E40").Value = "  +15.21%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0991"
$ws.Range("E41").Value = "  -14.14%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0223"
$ws.Range("E42").Value = "  -1.63%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.27"
$ws.Range("E43").Value = "  -2.70%  "

$ws.Range("E44").Value = "  -2.13%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "95.73"
$ws.Range("E45").Value = "  -3.17%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.45"
$ws.Range("E46").Value = "  -0.95%  "

$ws.Range("D47").Value = "1.275.68"
$ws.Range("E47").Value = "  -2.32%  "

$ws.Range("E48").Value = "  -3.68%  "

$ws.Range("E49").Value = "  -1.91%  "

$ws.Range("D50").Value = "2.231.26"
$ws.Range("E50").Value = "  -1.41%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.82"
$ws.Range("E51").Value = "  +0.14%  "
